# Weekly fruit/vegetable data update:
# Insert one new data row at row 343 (shifting the existing rows 343-400 down
# to 344-401) and populate the new row with this week's record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("343:343").Insert()

$ws.Cells.Item(343, 1).Value = 9
$ws.Cells.Item(343, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(343, 3).Value = "Metropolitana"
$ws.Cells.Item(343, 4).Value = 45218
$ws.Cells.Item(343, 5).Value = 13
$ws.Cells.Item(343, 6).Value = 100112001
$ws.Cells.Item(343, 7).Value = "Berenjena"
$ws.Cells.Item(343, 8).Value = "Sin especificar"
$ws.Cells.Item(343, 9).Value = "Primera"
$ws.Cells.Item(343, 10).Value = 70
$ws.Cells.Item(343, 11).Value = 10000
$ws.Cells.Item(343, 12).Value = 11000
$ws.Cells.Item(343, 13).Value = 10514
$ws.Cells.Item(343, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(343, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(343, 16).Value = 210
$ws.Cells.Item(343, 17).Value = 50
$ws.Cells.Item(343, 18).Value = "Hortaliza"
